$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New game rows to append (Away team, Away Pts, Home team, Home Pts, Overtime, Attend., Arena, Win, Loss)
$rows = @(
    @("Toronto Raptors", 123, "Atlanta Hawks", 121, "No", 17832, "State Farm Arena", "Toronto Raptors", "Atlanta Hawks"),
    @("Cleveland Cavaliers", 97, "Philadelphia 76ers", 104, "No", 17832, "Wells Fargo Center", "Philadelphia 76ers", "Cleveland Cavaliers"),
    @("Phoenix Suns", 110, "Houston Rockets", 114, "No", 17832, "Toyota Center", "Houston Rockets", "Phoenix Suns"),
    @("Los Angeles Clippers", 101, "Memphis Grizzlies", 95, "No", 17832, "FedEx Forum", "Los Angeles Clippers", "Memphis Grizzlies"),
    @("Miami Heat", 106, "New Orleans Pelicans", 95, "No", 17832, "Smoothie King Center", "Miami Heat", "New Orleans Pelicans"),
    @("Washington Wizards", 106, "Oklahoma City Thunder", 147, "No", 17832, "Paycom Center", "Oklahoma City Thunder", "Washington Wizards"),
    @("Charlotte Hornets", 84, "Golden State Warriors", 97, "No", 17832, "Chase Center", "Golden State Warriors", "Charlotte Hornets"),
    @("Milwaukee Bucks", 112, "Minnesota Timberwolves", 107, "No", 17832, "Target Center", "Milwaukee Bucks", "Minnesota Timberwolves"),
    @("Denver Nuggets", 127, "Portland Trail Blazers", 112, "No", 17832, "Moda Center", "Denver Nuggets", "Portland Trail Blazers"),
    @("San Antonio Spurs", 118, "Los Angeles Lakers", 123, "No", 17832, "Crypto.com Arena", "Los Angeles Lakers", "San Antonio Spurs")
)

$startRow = 835
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}

# Update sheet view to reflect scrolled position and new selection
$excel.ActiveWindow.ScrollRow = 812
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G835").Select()
